$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.848.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5019"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9096"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07652"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.909.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.503"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.83"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008725"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.886.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.176"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.585"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "153.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "

# Row 25
$ws.Range("E25").Value = "  -2.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.220"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.28%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.921"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09033"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.206"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.783"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.29%  "

# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.225"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.58%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7725"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02074"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.529"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.83%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.026"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05270"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.70%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.888"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.466"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1514"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.20%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4814"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47
$ws.Range("E47").Value = "  -0.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.632"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06062"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9063"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.79%  "
